$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidential disclosure text date from 2021-05-17 to 2021-05-18
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-35
$ws.Range("D2").Value = 0.03552155342395606
$ws.Range("E2").Value = -0.0003946329913182556
$ws.Range("D3").Value = 0.02014245647179362
$ws.Range("E3").Value = 0.001579155151993916
$ws.Range("D4").Value = 0.01923400593261784
$ws.Range("E4").Value = 0.00547451871366933
$ws.Range("D5").Value = 0.03762628893500599
$ws.Range("E5").Value = -0.001761183515322395
$ws.Range("D6").Value = 0.0342666131588425
$ws.Range("E6").Value = 0.0004001600640257674
$ws.Range("D7").Value = 0.0197691704909891
$ws.Range("E7").Value = -0.001918391026880739
$ws.Range("D8").Value = 0.03723760866685981
$ws.Range("E8").Value = -0.002546148949713722
$ws.Range("D9").Value = 0.02040059114430408
$ws.Range("E9").Value = -8.99523252676504E-05
$ws.Range("D10").Value = 0.0261116678502138
$ws.Range("E10").Value = 0.0003904343582235015
$ws.Range("D11").Value = 0.02422561277855366
$ws.Range("E11").Value = -0.008153603366649165
$ws.Range("D12").Value = 0.05762550717027623
$ws.Range("E12").Value = -0.00590458195559751
$ws.Range("D13").Value = 0.02462678178974191
$ws.Range("E13").Value = 0.002235469448584215
$ws.Range("D14").Value = 0.02693531319745393
$ws.Range("E14").Value = -0.004345127250155234
$ws.Range("D15").Value = 0.03307875528492639
$ws.Range("E15").Value = -0.0131170176044183
$ws.Range("D16").Value = 0.01983248077204321
$ws.Range("E16").Value = -0.001495886312640304
$ws.Range("D17").Value = 0.03137478098937354
$ws.Range("E17").Value = -0.004630381803411909
$ws.Range("D18").Value = 0.04184952306168283
$ws.Range("E18").Value = -0.0004628558204118338
$ws.Range("D19").Value = 0.1253826982876652
$ws.Range("E19").Value = 0.0006675567423231055
$ws.Range("D20").Value = 0.00936405953294705
$ws.Range("E20").Value = -0.006134969325153561
$ws.Range("D21").Value = 0.0152905033382558
$ws.Range("E21").Value = -0.0002100252030242666
$ws.Range("D22").Value = 0.01718140098374859
$ws.Range("E22").Value = 0.003681847985070874
$ws.Range("D23").Value = 0.01536293804789983
$ws.Range("E23").Value = 0.003616636528029016
$ws.Range("D24").Value = 0.02135350772564502
$ws.Range("E24").Value = 0.001334428248819464
$ws.Range("D25").Value = 0.01281757929254097
$ws.Range("E25").Value = -0.009059419131361657
$ws.Range("D26").Value = 0.04243646849742961
$ws.Range("E26").Value = -0.003086469802865865
$ws.Range("D27").Value = 0.02390783798620962
$ws.Range("E27").Value = 0.0001961553550411388
$ws.Range("D28").Value = 0.04542803220062028
$ws.Range("E28").Value = 0.001434034416825769
$ws.Range("D29").Value = 0.05571977595651714
$ws.Range("E29").Value = -0.01004124081047164
$ws.Range("D30").Value = 0.01310726715665569
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0.02059011419177523
$ws.Range("E31").Value = 0.0007674597083653367
$ws.Range("D32").Value = 0.01338716792418233
$ws.Range("E32").Value = -0.009519276534983301
$ws.Range("D33").Value = 0.04185222470813683
$ws.Range("E33").Value = -0.002575991756826435
$ws.Range("D34").Value = 0.01695971305113643
$ws.Range("E34").Value = 0.0008836524300439574
$ws.Range("E35").Value = -0.002052681410561252
